$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (shifting existing batsman..sr columns right to F..K)
$ws.Range("D1:E1").EntireColumn.Insert()

# Set header values for the new columns
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Set data values for the new columns
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"
